$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "process" column in Y: header + same value for every data row
$ws.Range("Y1").Value = "process"
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 25).Value = "pp->Z/gamma*->l+ l-"
}

# Update the view: scroll so column K is at the left edge, select Y1:Y35
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("Y1:Y35").Select()
